$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten a few display names down to the short form used for avatar/photo lookups
$ws.Range("B2").Value = "Trung"
$ws.Range("E4").Value = "Việt Trì Phú Thọ"
$ws.Range("B5").Value = "Khương"
$ws.Range("B7").Value = "An"

# Swap the old last record (id 46 / M.C.Duy) for a new employee (id 28 / T.H.Giang)
$ws.Range("A14").Value = 28
$ws.Range("B14").Value = "T.H.Giang"
$ws.Range("C14").Value = "giang"

# The old trailing record (id 47 / N.Đ.Manh) is no longer needed
$ws.Rows.Item(15).Delete()
